$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (F:S) per revised ERA operativo Abril-Diciembre 2025 ---
# Row 2
$ws.Range("F2").Value = 5.733927907840951
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("S2").Value = 0

# Row 3
$ws.Range("F3").Value = 5.881428122230798
$ws.Range("G3").Value = 9.803921568627452
$ws.Range("H3").Value = 1
$ws.Range("L3").Value = 0
$ws.Range("O3").Value = 0.1111111111111111
$ws.Range("S3").Value = 0

# Row 4
$ws.Range("F4").Value = 5.881428122230798
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("S4").Value = 0

# Row 5
$ws.Range("F5").Value = 5.934025097904697
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("S5").Value = 0

# Row 6
$ws.Range("F6").Value = 5.934025097904697
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("S6").Value = 0

# Row 7
$ws.Range("F7").Value = 5.994625961181145
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0

# Row 8
$ws.Range("F8").Value = 5.994625961181145
$ws.Range("G8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0

# Row 9
$ws.Range("F9").Value = 5.994625961181145
$ws.Range("G9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0

# Row 10
$ws.Range("F10").Value = 6.247320126918789
$ws.Range("G10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0

# Row 11
$ws.Range("F11").Value = 6.334219478032186
$ws.Range("G11").Value = 11.76470588235294
$ws.Range("I11").Value = 1
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("P11").Value = 0.1111111111111111
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0

# Row 12
$ws.Range("F12").Value = 6.35508675642455
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("R12").Value = 0
$ws.Range("S12").Value = 0

# Row 13
$ws.Range("F13").Value = 6.542892261955807
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 0

# Row 14
$ws.Range("F14").Value = 6.734413858159677
$ws.Range("G14").Value = 78.43137254901961
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0
$ws.Range("R14").Value = 0.2222222222222222
$ws.Range("S14").Value = 0

# Row 15
$ws.Range("F15").Value = 6.812451762284538
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("R15").Value = 0
$ws.Range("S15").Value = 0

# Row 16
$ws.Range("F16").Value = 6.812451762284538
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0
$ws.Range("R16").Value = 0
$ws.Range("S16").Value = 0

# --- Append new row 17: Recluta1 ---
$ws.Range("A17").Value = "Recluta1"
$ws.Range("B17").Value = 44713
$ws.Range("B17").NumberFormat = $ws.Range("B16").NumberFormat
$ws.Range("C17").Value = 1.706849315068493
$ws.Range("D17").Value = "Padawan-Sin Fijo"
$ws.Range("E17").Value = 15
$ws.Range("F17").Value = 6.812451762284538
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("R17").Value = 0
$ws.Range("S17").Value = 0
